$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.642.19"
$ws.Cells.Item(2, 5).Value = "  -1.26%  "
$ws.Cells.Item(3, 4).Value = "2.300.87"
$ws.Cells.Item(3, 5).Value = "  -0.17%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "302.16"
$ws.Cells.Item(5, 5).Value = "  -2.12%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "98.64"
$ws.Cells.Item(6, 5).Value = "  -5.30%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.498"
$ws.Cells.Item(7, 5).Value = "  -5.10%  "
$ws.Cells.Item(8, 5).Value = "  +0.17%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.500"
$ws.Cells.Item(9, 5).Value = "  -3.32%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "34.54"
$ws.Cells.Item(10, 5).Value = "  -3.51%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0788"
$ws.Cells.Item(11, 5).Value = "  -2.33%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.112"
$ws.Cells.Item(12, 5).Value = "  +0.30%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.68"
$ws.Cells.Item(13, 5).Value = "  -3.81%  "
$ws.Cells.Item(14, 4).Value = "2.660.11"
$ws.Cells.Item(14, 5).Value = "  +0.03%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.56"
$ws.Cells.Item(15, 5).Value = "  +3.15%  "
$ws.Cells.Item(16, 4).Value = "2.304.88"
$ws.Cells.Item(16, 5).Value = "  +0.28%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.795"
$ws.Cells.Item(17, 5).Value = "  -0.33%  "
$ws.Cells.Item(18, 4).Value = "42.562.51"
$ws.Cells.Item(18, 5).Value = "  -1.35%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0902"
$ws.Cells.Item(19, 5).Value = "  -2.00%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.53"
$ws.Cells.Item(20, 5).Value = "  -3.62%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.03"
$ws.Cells.Item(21, 5).Value = "  -2.38%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "67.64"
$ws.Cells.Item(22, 5).Value = "  -0.17%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "234.69"
$ws.Cells.Item(23, 5).Value = "  -2.31%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.94"
$ws.Cells.Item(24, 5).Value = "  -3.33%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.50"
$ws.Cells.Item(25, 5).Value = "  -3.51%  "
$ws.Cells.Item(26, 5).Value = "  +0.01%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "24.62"
$ws.Cells.Item(27, 5).Value = "  -0.29%  "
$ws.Cells.Item(28, 5).Value = "  -2.65%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "34.26"
$ws.Cells.Item(29, 5).Value = "  -4.87%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "163.38"
$ws.Cells.Item(30, 5).Value = "  +0.84%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "9.09"
$ws.Cells.Item(31, 5).Value = "  -4.79%  "
$ws.Cells.Item(32, 5).Value = "  +0.06%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.98"
$ws.Cells.Item(33, 5).Value = "  -4.87%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.42"
$ws.Cells.Item(34, 5).Value = "  -4.71%  "
# Rows 35-38: coins reordered (Celestia/RenderToken swap, Hedera/LidoDAOToken swap)
# along with refreshed price/volume figures.
$ws.Cells.Item(35, 2).Value = "RenderToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.43"
$ws.Cells.Item(35, 5).Value = "  -1.35%  "

$ws.Cells.Item(36, 2).Value = "Celestia"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "16.64"
$ws.Cells.Item(36, 5).Value = "  -8.98%  "

$ws.Cells.Item(37, 2).Value = "LidoDAOToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.87"
$ws.Cells.Item(37, 5).Value = "  -3.98%  "

$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0690"
$ws.Cells.Item(38, 5).Value = "  -5.62%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.79"
$ws.Cells.Item(39, 5).Value = "  -3.44%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0998"
$ws.Cells.Item(40, 5).Value = "  -5.69%  "
$ws.Cells.Item(41, 5).Value = "  -4.38%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.51"
$ws.Cells.Item(42, 5).Value = "  +0.76%  "
$ws.Cells.Item(43, 4).Value = "1.959.22"
$ws.Cells.Item(43, 5).Value = "  -0.21%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0278"
$ws.Cells.Item(44, 5).Value = "  -3.28%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "18.38"
$ws.Cells.Item(45, 5).Value = "  -1.36%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.12"
$ws.Cells.Item(46, 5).Value = "  -0.47%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.86"
$ws.Cells.Item(47, 5).Value = "  -6.25%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "55.07"
$ws.Cells.Item(48, 5).Value = "  -3.79%  "
$ws.Cells.Item(49, 4).Value = "2.527.98"
$ws.Cells.Item(49, 5).Value = "  -0.02%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.82"
$ws.Cells.Item(50, 5).Value = "  -3.99%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.66"
$ws.Cells.Item(51, 5).Value = "  +0.48%  "
